$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.080.41"
$ws.Range("E2").Value = "  -1.05%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.877.32"
$ws.Range("E3").Value = "  -1.60%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.61%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.13"
$ws.Range("E5").Value = "  -0.91%  "

# Row 6
$ws.Range("E6").Value = "  +0.57%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4646"
$ws.Range("E7").Value = "  -2.44%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2784"
$ws.Range("E8").Value = "  -2.64%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06532"
$ws.Range("E9").Value = "  -2.50%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.00"
$ws.Range("E10").Value = "  +1.13%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07741"
$ws.Range("E11").Value = "  +0.28%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.39"
$ws.Range("E12").Value = "  -5.12%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.896.40"
$ws.Range("E13").Value = "  -0.64%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.089"
$ws.Range("E14").Value = "  -2.18%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6525"
$ws.Range("E15").Value = "  -3.52%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "279.19"
$ws.Range("E16").Value = "  +7.38%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.114.63"
$ws.Range("E17").Value = "  -1.00%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.004"
$ws.Range("E18").Value = "  +0.39%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.142.93"
$ws.Range("E19").Value = "  -0.56%  "

# Row 20
$ws.Range("B20").Value = "BinanceUSD"
$ws.Range("C20").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.006"
$ws.Range("E20").Value = "  +0.57%  "

# Row 21
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.35"
$ws.Range("E21").Value = "  -2.99%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.292"
$ws.Range("E22").Value = "  -2.22%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.000007208"
$ws.Range("E23").Value = "  -3.83%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.135"
$ws.Range("E24").Value = "  -2.83%  "

# Row 25
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "164.95"
$ws.Range("E25").Value = "  +0.29%  "

# Row 26
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.155"
$ws.Range("E26").Value = "  -3.19%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.79"
$ws.Range("E27").Value = "  -1.14%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.985"
$ws.Range("E28").Value = "  -3.46%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.380"
$ws.Range("E29").Value = "  +0.24%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09786"
$ws.Range("E30").Value = "  -3.15%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.431"
$ws.Range("E31").Value = "  -4.18%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.485"
$ws.Range("E32").Value = "  -1.83%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.144"
$ws.Range("E33").Value = "  -2.63%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04636"
$ws.Range("E34").Value = "  -3.15%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7033"
$ws.Range("E35").Value = "  -3.83%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.082"
$ws.Range("E36").Value = "  -2.90%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.710"
$ws.Range("E37").Value = "  +0.08%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01843"
$ws.Range("E38").Value = "  -4.36%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.664"
$ws.Range("E39").Value = "  +6.86%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.517"
$ws.Range("E40").Value = "  -3.27%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "71.83"
$ws.Range("E41").Value = "  -4.07%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8652"
$ws.Range("E42").Value = "  -0.13%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.921"
$ws.Range("E43").Value = "  -3.44%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.005"
$ws.Range("E44").Value = "  +0.57%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "103.33"
$ws.Range("E45").Value = "  -3.13%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4132"
$ws.Range("E46").Value = "  -3.01%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "970.42"
$ws.Range("E47").Value = "  -7.43%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.143"
$ws.Range("E48").Value = "  -4.59%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.043"
$ws.Range("E49").Value = "  +1.64%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1155"
$ws.Range("E50").Value = "  -3.72%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05745"
$ws.Range("E51").Value = "  +0.67%  "
